$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185, shifting existing rows 185:248 down to 186:249
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new record
$ws.Cells.Item(185, 1).Value = 11
$ws.Cells.Item(185, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(185, 3).Value = "Bíobío"
$ws.Cells.Item(185, 4).Value = 44988
$ws.Cells.Item(185, 5).Value = 8
$ws.Cells.Item(185, 6).Value = "Fruta"
$ws.Cells.Item(185, 7).Value = 100108
$ws.Cells.Item(185, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(185, 9).Value = 100108005
$ws.Cells.Item(185, 10).Value = "Piña"
$ws.Cells.Item(185, 11).Value = "Caramelo"
$ws.Cells.Item(185, 12).Value = "Primera"
$ws.Cells.Item(185, 13).Value = 100
$ws.Cells.Item(185, 14).Value = 22000
$ws.Cells.Item(185, 15).Value = 23000
$ws.Cells.Item(185, 16).Value = 22500
$ws.Cells.Item(185, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(185, 18).Value = "Ecuador"
$ws.Cells.Item(185, 19).Value = 1875
$ws.Cells.Item(185, 20).Value = 12
